$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (N3:O3) - blank cells, same border-only style as M3 ---
$ws.Range("M3").Copy() | Out-Null
$ws.Range("N3:O3").PasteSpecial(-4122) | Out-Null

# --- Row 4 (N4:O4) - year headers, same style as M4 ---
$ws.Range("N4").Value = 2021
$ws.Range("O4").Value = 2022
$ws.Range("M4").Copy() | Out-Null
$ws.Range("N4:O4").PasteSpecial(-4122) | Out-Null

# --- Row 5 (N5:O5) - new style (bold, numFmt 0.0, no border) ---
$ws.Range("N5").Value = 40.007977647471066
$ws.Range("O5").Value = 42.620582506455563
$ws.Range("N5").NumberFormat = "0.0"
$ws.Range("N5").Font.Name = "Times New Roman"
$ws.Range("N5").Font.Size = 10
$ws.Range("N5").Font.Bold = $true
$ws.Range("N5").Copy() | Out-Null
$ws.Range("O5").PasteSpecial(-4122) | Out-Null

# --- Rows 6-13 (N:O) - new style (regular, numFmt 0.0, no border) ---
$ws.Range("N6").Value = 5.7072514621689896
$ws.Range("O6").Value = 8.1443914479075037
$ws.Range("N7").Value = 8.9893229854028949
$ws.Range("O7").Value = 10.715961386284755
$ws.Range("N8").Value = 66.307512472824584
$ws.Range("O8").Value = 81.977461999426666
$ws.Range("N9").Value = 23.475213049310256
$ws.Range("O9").Value = 29.828871240443185
$ws.Range("N10").Value = 9.8045372040896162
$ws.Range("O10").Value = 9.7218425128664112
$ws.Range("N11").Value = 9.3737779268960448
$ws.Range("O11").Value = 8.6167819403064012
$ws.Range("N12").Value = 70.457032471318783
$ws.Range("O12").Value = 69.915337594090886
$ws.Range("N13").Value = 98.411252120183207
$ws.Range("O13").Value = 99.08571752721997

$ws.Range("N6").NumberFormat = "0.0"
$ws.Range("N6").Font.Name = "Times New Roman"
$ws.Range("N6").Font.Size = 10
$ws.Range("N6").Font.Bold = $false
$ws.Range("N6").Copy() | Out-Null
$ws.Range("O6").PasteSpecial(-4122) | Out-Null
$ws.Range("N7:O13").PasteSpecial(-4122) | Out-Null

# --- Row 14 (N14:O14) - new style (regular, numFmt 0.0, bottom medium border) ---
$ws.Range("N14").Value = 63.900563564170795
$ws.Range("O14").Value = 64.805252627098838
$ws.Range("M3").Copy() | Out-Null
$ws.Range("N14").PasteSpecial(-4122) | Out-Null
$ws.Range("N14").NumberFormat = "0.0"
$ws.Range("N14").Copy() | Out-Null
$ws.Range("O14").PasteSpecial(-4122) | Out-Null

# --- Update selection to match target ---
$ws.Range("P8").Select() | Out-Null
